$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 63 (pushes existing rows 63-66 down to 64-67),
# matching a newly reported weekly price observation being prepended to the
# (reverse-chronological) data table.
$ws.Rows.Item(63).Insert()

# Populate the new row 63 with the latest observation.
$ws.Cells.Item(63, 1).Value  = 7
$ws.Cells.Item(63, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(63, 3).Value  = "Ñuble"
$ws.Cells.Item(63, 4).Value  = 44595
$ws.Cells.Item(63, 5).Value  = 16
$ws.Cells.Item(63, 6).Value  = 100112031
$ws.Cells.Item(63, 7).Value  = "Poroto verde"
$ws.Cells.Item(63, 8).Value  = "Sin especificar"
$ws.Cells.Item(63, 9).Value  = "Primera"
$ws.Cells.Item(63, 10).Value = 60
$ws.Cells.Item(63, 11).Value = 29000
$ws.Cells.Item(63, 12).Value = 30000
$ws.Cells.Item(63, 13).Value = 29500
$ws.Cells.Item(63, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(63, 15).Value = "Región del Maule"
$ws.Cells.Item(63, 16).Value = 1180
$ws.Cells.Item(63, 17).Value = 25
$ws.Cells.Item(63, 18).Value = "Hortaliza"
